$wb = $excel.ActiveWorkbook

# --- Sheet "FBS" (sheet1) updates ---
$fbs = $wb.Worksheets.Item("FBS")

# Timestamp column (AK2:AK9) - update to new timestamp value
$fbs.Range("AK2:AK9").Value = "2024-12-04T10:01:12.142724"

# temp_fg column (O)
$fbs.Range("O2").Value = 27.02
$fbs.Range("O3").Value = 35.6
$fbs.Range("O4").Value = 38.3
$fbs.Range("O5").Value = 66.2
$fbs.Range("O6").Value = 37.94
$fbs.Range("O7").Value = 55.7
$fbs.Range("O8").Value = 44.54
$fbs.Range("O9").Value = 45.44

# gs_fg / away_fg columns (S2, T2)
$fbs.Range("S2").Value = -0.37
$fbs.Range("T2").Value = -0.62

# Odds_n column (Z5)
$fbs.Range("Z5").Value = -106

# wind_dir_fg column (Q8): NE -> N
$fbs.Range("Q8").Value = "N"

# --- Sheet "Other" (sheet2) updates ---
$other = $wb.Worksheets.Item("Other")

# temp_fg column (Q)
$other.Range("Q2").Value = 55.58000000000001
$other.Range("Q3").Value = 43.04
$other.Range("Q5").Value = 60.8
